$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-85 (inclusive) in column C get corrected to the flat-lined
# fitness value of 7293 (matches rows 86 onward already in the sheet).
$ws.Range("C2:C85").Value = 7293
